$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 22083.25
$ws.Range("I69").Value = 9999
$ws.Range("K69").Value = 29997
$ws.Range("M69").Value = -29123
$ws.Range("H72").Value = 22083.25
$ws.Range("I72").Value = 9999
$ws.Range("K72").Value = 89991
$ws.Range("M72").Value = -85623
$ws.Range("H74").Value = 1994.5
$ws.Range("I74").Value = 1994.5
$ws.Range("K74").Value = 1994.5
$ws.Range("M74").Value = -1058.5
$ws.Range("H77").Value = 1994.5
$ws.Range("I77").Value = 1994.5
$ws.Range("K77").Value = 9972.5
$ws.Range("M77").Value = -5292.5
$ws.Range("H137").Value = 2115.682
$ws.Range("I137").Value = 1585.4117
$ws.Range("K137").Value = 4756.2351
$ws.Range("M137").Value = -2206.2351
$ws.Range("H141").Value = 5801.091
$ws.Range("I141").Value = 4358.625
$ws.Range("K141").Value = 13075.875
$ws.Range("M141").Value = -7895.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 536.7778
$ws.Range("I2").Value = 334.1875
$ws.Range("K2").Value = 334.1875
$ws.Range("M2").Value = -221.1875
$ws.Range("H45").Value = 9787.5
$ws.Range("I45").Value = 12827.75
$ws.Range("K45").Value = 12827.75
$ws.Range("M45").Value = -12450.75
$ws.Range("H102").Value = 7635.3335
$ws.Range("J102").Value = 2769.5
$ws.Range("L102").Value = 2769.5
$ws.Range("N102").Value = -6013.5
$ws.Range("H110").Value = 72447.42999999999
$ws.Range("I110").Value = 77961.62
$ws.Range("K110").Value = 77961.62
$ws.Range("M110").Value = -75916.62
$ws.Range("H116").Value = 536.7778
$ws.Range("I116").Value = 334.1875
$ws.Range("K116").Value = 334.1875
$ws.Range("M116").Value = 1959.8125
$ws.Range("H122").Value = 4204.278
$ws.Range("I122").Value = 1585.2
$ws.Range("J122").Value = 17299.666
$ws.Range("K122").Value = 4755.6
$ws.Range("L122").Value = 51898.99800000001
$ws.Range("M122").Value = -2305.6
$ws.Range("N122").Value = -56798.99800000001
$ws.Range("H132").Value = 2085070.1
$ws.Range("I132").Value = 2440563
$ws.Range("K132").Value = 7321689
$ws.Range("M132").Value = -7319159

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 536.7778
$ws.Range("I3").Value = 334.1875
$ws.Range("K3").Value = 334.1875
$ws.Range("M3").Value = -220.1875
$ws.Range("H94").Value = 3548.5
$ws.Range("I94").Value = 3548.5
$ws.Range("K94").Value = 3548.5
$ws.Range("M94").Value = -3097.5
$ws.Range("H105").Value = 3178.4375
$ws.Range("I105").Value = 2832.2727
$ws.Range("J105").Value = 3940
$ws.Range("K105").Value = 2832.2727
$ws.Range("L105").Value = 3940
$ws.Range("M105").Value = -1085.2727
$ws.Range("N105").Value = -7434
$ws.Range("H134").Value = 29413264
$ws.Range("I134").Value = 29413264
$ws.Range("K134").Value = 88239792
$ws.Range("M134").Value = -88237257

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 100000
$ws.Range("I16").Value = 100000
$ws.Range("K16").Value = 100000
$ws.Range("M16").Value = -99713
$ws.Range("H18").Value = 67842
$ws.Range("J18").Value = 70684
$ws.Range("L18").Value = 70684
$ws.Range("N18").Value = -71144
$ws.Range("H99").Value = 1795.1428
$ws.Range("I99").Value = 1909.8182
$ws.Range("J99").Value = 1374.6666
$ws.Range("K99").Value = 1909.8182
$ws.Range("L99").Value = 1374.6666
$ws.Range("M99").Value = -411.8181999999999
$ws.Range("N99").Value = -4370.6666
$ws.Range("H105").Value = 3402343.5
$ws.Range("I105").Value = 4082432.5
$ws.Range("K105").Value = 4082432.5
$ws.Range("M105").Value = -4080685.5
$ws.Range("H113").Value = 100000
$ws.Range("I113").Value = 100000
$ws.Range("K113").Value = 100000
$ws.Range("M113").Value = -97830
$ws.Range("H126").Value = 1795.1428
$ws.Range("I126").Value = 1909.8182
$ws.Range("J126").Value = 1374.6666
$ws.Range("K126").Value = 5729.4546
$ws.Range("L126").Value = 4123.9998
$ws.Range("M126").Value = -3259.4546
$ws.Range("N126").Value = -9063.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 10950
$ws.Range("I57").Value = 1900
$ws.Range("K57").Value = 5700
$ws.Range("M57").Value = -5141
$ws.Range("H68").Value = 2293.85
$ws.Range("I68").Value = 1145
$ws.Range("J68").Value = 2676.8
$ws.Range("K68").Value = 3435
$ws.Range("L68").Value = 8030.400000000001
$ws.Range("M68").Value = -2624
$ws.Range("N68").Value = -9652.400000000001
$ws.Range("H71").Value = 2293.85
$ws.Range("I71").Value = 1145
$ws.Range("J71").Value = 2676.8
$ws.Range("K71").Value = 10305
$ws.Range("L71").Value = 24091.2
$ws.Range("M71").Value = -6249
$ws.Range("N71").Value = -32203.2
$ws.Range("H75").Value = 3880
$ws.Range("J75").Value = 3880
$ws.Range("L75").Value = 11640
$ws.Range("N75").Value = -13636
$ws.Range("H78").Value = 3880
$ws.Range("J78").Value = 3880
$ws.Range("L78").Value = 34920
$ws.Range("N78").Value = -44904
$ws.Range("H107").Value = 975.2
$ws.Range("I107").Value = 317
$ws.Range("J107").Value = 1492.3572
$ws.Range("K107").Value = 951
$ws.Range("L107").Value = 4477.071599999999
$ws.Range("M107").Value = 969
$ws.Range("N107").Value = -8317.071599999999
$ws.Range("H122").Value = 703.7143
$ws.Range("I122").Value = 648.5
$ws.Range("J122").Value = 745.125
$ws.Range("K122").Value = 5836.5
$ws.Range("L122").Value = 6706.125
$ws.Range("M122").Value = -3386.5
$ws.Range("N122").Value = -11606.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 6962.625
$ws.Range("I19").Value = 7516.8335
$ws.Range("J19").Value = 5300
$ws.Range("K19").Value = 7516.8335
$ws.Range("L19").Value = 5300
$ws.Range("M19").Value = -7228.8335
$ws.Range("N19").Value = -5876
$ws.Range("H97").Value = 1083.3478
$ws.Range("I97").Value = 933.2105
$ws.Range("J97").Value = 1796.5
$ws.Range("K97").Value = 933.2105
$ws.Range("L97").Value = 1796.5
$ws.Range("M97").Value = -437.2105
$ws.Range("N97").Value = -2788.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1632.3572
$ws.Range("I22").Value = 1962.7778
$ws.Range("K22").Value = 1962.7778
$ws.Range("M22").Value = -1667.7778
$ws.Range("H27").Value = 1632.3572
$ws.Range("I27").Value = 1962.7778
$ws.Range("K27").Value = 1962.7778
$ws.Range("M27").Value = -1855.7778
$ws.Range("H61").Value = 2964.75
$ws.Range("I61").Value = 3034.2727
$ws.Range("K61").Value = 3034.2727
$ws.Range("M61").Value = -2832.2727
$ws.Range("H74").Value = 80894.25
$ws.Range("I74").Value = 89580
$ws.Range("J74").Value = 77999
$ws.Range("K74").Value = 89580
$ws.Range("L74").Value = 77999
$ws.Range("M74").Value = -88582
$ws.Range("N74").Value = -79995
$ws.Range("H77").Value = 80894.25
$ws.Range("I77").Value = 89580
$ws.Range("J77").Value = 77999
$ws.Range("K77").Value = 268740
$ws.Range("L77").Value = 233997
$ws.Range("M77").Value = -263748
$ws.Range("N77").Value = -243981
$ws.Range("H113").Value = 2964.75
$ws.Range("I113").Value = 3034.2727
$ws.Range("K113").Value = 3034.2727
$ws.Range("M113").Value = -864.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2412.25
$ws.Range("I81").Value = 1899.7142
$ws.Range("J81").Value = 6000
$ws.Range("K81").Value = 3799.4284
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -2738.4284
$ws.Range("N81").Value = -14122
$ws.Range("H84").Value = 2412.25
$ws.Range("I84").Value = 1899.7142
$ws.Range("J84").Value = 6000
$ws.Range("K84").Value = 18997.142
$ws.Range("L84").Value = 60000
$ws.Range("M84").Value = -13693.142
$ws.Range("N84").Value = -70608

Write-Host "Applied 202 cell updates across 8 sheets"